# "Generate Report for Archive"
# 1. Update the shared "Status" text from "Ready for handoff" to "In Translation"
#    wherever it appears (Overview!E2:F2, zh-cn!C2, de-de!C2).
# 2. Shrink the "Status" column (Overview columns E & F; zh-cn/de-de column C)
#    from their old width down to the new, narrower width.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- 1. Replace the status text everywhere it is used ---
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# --- 2. Narrow the Status columns ---
# Target stored width is 13.4101845877511 characters; the engine stores
# ColumnWidth internally as round(input*6)/6 + 5/6, so 12.5 is the closest
# input that reproduces that value.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
